# Amplitude Estimation Benchmark-Results.xlsx edit
# - Center-align (style s="1") the second results block's header/table rows (19:27),
#   matching the alignment already used by the rest of the sheet.
# - Append a third results block (rows 28:40): blank separator row, 3 preamble
#   lines (merged A:S), a header row, 7 data rows, and a trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Apply the "centered" style (s="1") to A19:T27, matching rows above them.
# ---------------------------------------------------------------------------
$ws.Range("A19:T27").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2) New block starting at row 28 (row 28 left blank, like row 13/14/27).
# ---------------------------------------------------------------------------

# Row 28: blank separator row, centered style like its counterparts.
$ws.Range("A28:T28").HorizontalAlignment = -4108

# Rows 29-31: preamble text lines, merged A:S, centered style.
$ws.Cells.Item(29, 1).Value = "Qsim: Algorithm = Amplitude Estimation Simulator = dm_simulator"
$ws.Cells.Item(30, 1).Value = "CPU: 12th Gen Intel(R) Core(TM) i9-12900 with 24 cores"
$ws.Cells.Item(31, 1).Value = "Configuration: Min_Qubits = 4 Max_Qubits = 10 Skip_Qubits = 1 num_circuits = 2  QV_ = None Last_Updated = 2025-02-10 11:55:00"

$ws.Range("A29:T31").HorizontalAlignment = -4108

$ws.Range("A29:S29").Merge()
$ws.Range("A30:S30").Merge()
$ws.Range("A31:S31").Merge()

# Row 32: column headers (not centered, matches row 19's pre-edit state).
$headers = @(
    "Number of Qubits",
    "avg_creation_times (ms)",
    "std_creation_times (ms)",
    "avg_elapsed_times (ms)",
    "std_elapsed_times (ms)",
    "avg_quantum_times (ms)",
    "std_quantum_times (ms)",
    "avg_circuit_depths",
    "avg_transpiled_depths",
    "Average_Rescaled_fidelity",
    "Average_Hellinger_fidelity",
    "std_Rescaled_Fidelity",
    "std_hellinger_fidelity",
    "avg_1Q_algorithmic_gate_counts",
    "avg_2Q_algorithmic_gate_counts",
    "avg_xi (n2q/n1q+n2q)",
    "avg_1Q_Transpiled_gate_counts",
    "avg_2Q_Transpiled_gate_counts",
    "avg_tr_xi (tr_n2q/tr_n1q+tr_n2q)",
    "max_memory (MB)"
)
$col = 1
foreach ($h in $headers) {
    $ws.Cells.Item(32, $col).Value = $h
    $col++
}

# Rows 33-39: data rows for qubit counts 4..10.
$data = @(
    @(4, 110.604, 28.004, 2543.473, 643.704, 14.002, 0.082, 208, 167.5, 1, 1, 0, 0, 37, 62, 0.63, 128.5, 85, 0.4, 169.5),
    @(5, 7.714, 0.045, 326.915, 0.495, 50.646, 0.016, 479, 363, 1, 1, 0, 0, 79, 144, 0.65, 277, 188, 0.4, 170.84),
    @(6, 12.003, 0.107, 984.596, 7.997, 461.558, 0.404, 1017, 836, 1, 1, 0, 0, 161, 307, 0.66, 621, 426, 0.41, 173.34),
    @(7, 21.329, 0.1, 4629.913, 15.501, 3570.347, 9.883, 2086, 1701, 0.99, 0.99, 0.001, 0.001, 323, 631, 0.66, 1258, 868, 0.41, 177.37),
    @(8, 68.711, 0.194, 31408.188, 151.431, 29239.596, 153.081, 4214, 3417, 0.96, 0.97, 0.004, 0.004, 645, 1276, 0.66, 2521, 1744, 0.41, 189.78),
    @(9, 108.442, 0.103, 243110.663, 122.646, 238635.866, 69.383, 8457, 6832, 0.96, 0.96, 0.001, 0.001, 1287, 2562, 0.67, 5034, 3486, 0.41, 231.02),
    @(10, 188.401, 1.755, 2093656.623, 905.23, 2084767.013, 876.811, 16927, 13642, 0.9, 0.9, 0.003, 0.003, 2569, 5129, 0.67, 10045, 6958, 0.41, 382.17)
)

$r = 33
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Row 40: trailing blank row (default/general style, no value) - force the
# cells into existence (so the used range/dimension extends to row 40)
# without leaving any alignment/number-format behind.
$ws.Range("A40:T40").HorizontalAlignment = -4108
$ws.Range("A40:T40").ClearFormats()

Write-Host "Edit complete"
